$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for L1:O1
$ws.Range("L1").Value = "# Trues"
$ws.Range("M1").Value = "# Lies"
$ws.Range("N1").Value = "# Free Trues"
$ws.Range("O1").Value = "# Free Lies"

# Updated/new data values per row
# Row 2
$ws.Range("B2").Value = -0.04762177186784809
$ws.Range("C2").Value = 0.04495756990159295
$ws.Range("D2").Value = 0.1177978128723235
$ws.Range("E2").Value = 0.0479405232405964
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 2

# Row 3
$ws.Range("B3").Value = 0.107481103508598
$ws.Range("C3").Value = 0.1617039681492539
$ws.Range("D3").Value = -0.0381621483487045
$ws.Range("E3").Value = 0.2485940006675374
$ws.Range("L3").Value = 10
$ws.Range("M3").Value = 10
$ws.Range("N3").Value = 4
$ws.Range("O3").Value = 6

# Row 4
$ws.Range("B4").Value = 0.08311790308958619
$ws.Range("C4").Value = 0.06189214622931673
$ws.Range("D4").Value = 0.0511929751519607
$ws.Range("E4").Value = 0.07789183443763835
$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 10
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 4

# Row 5
$ws.Range("B5").Value = 0.01299609428079811
$ws.Range("C5").Value = 0.1275747575143158
$ws.Range("D5").Value = 0.2123092255702406
$ws.Range("E5").Value = 0.1087258918827383
$ws.Range("L5").Value = 10
$ws.Range("M5").Value = 10
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 3

# Row 6
$ws.Range("B6").Value = 0.2717126667191195
$ws.Range("C6").Value = -0.1275451596305487
$ws.Range("D6").Value = 0.2140809321972516
$ws.Range("E6").Value = 0.8197452730118098
$ws.Range("L6").Value = 10
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 7
$ws.Range("O6").Value = 3

# Row 7
$ws.Range("B7").Value = 0.09330843560988983
$ws.Range("C7").Value = 0.06872658849864427
$ws.Range("D7").Value = -0.01745157640583462
$ws.Range("E7").Value = -0.07159013563472139
$ws.Range("L7").Value = 10
$ws.Range("M7").Value = 10
$ws.Range("N7").Value = 8
$ws.Range("O7").Value = 2

# Row 8
$ws.Range("B8").Value = -0.1113193579371374
$ws.Range("C8").Value = 0.144591303632018
$ws.Range("D8").Value = 0.1677156367630634
$ws.Range("E8").Value = -0.1054721431556653
$ws.Range("L8").Value = 10
$ws.Range("M8").Value = 10
$ws.Range("N8").Value = 10
$ws.Range("O8").Value = 2

# Row 9
$ws.Range("B9").Value = -0.1019678352974485
$ws.Range("C9").Value = 0.007676501988482948
$ws.Range("D9").Value = -0.05500193102709847
$ws.Range("E9").Value = -0.2480356728751577
$ws.Range("L9").Value = 10
$ws.Range("M9").Value = 10
$ws.Range("N9").Value = 8
$ws.Range("O9").Value = 2

# Row 10
$ws.Range("B10").Value = 0.04163922900789983
$ws.Range("C10").Value = 0.200278241297822
$ws.Range("D10").Value = -0.1865226107058604
$ws.Range("E10").Value = 0.3053562308869803
$ws.Range("L10").Value = 10
$ws.Range("M10").Value = 10
$ws.Range("N10").Value = 8
$ws.Range("O10").Value = 1

# Row 11
$ws.Range("B11").Value = 0.06149143931986877
$ws.Range("C11").Value = 0.009947458025988945
$ws.Range("D11").Value = -0.103795218466967
$ws.Range("E11").Value = 0.07483201451396219
$ws.Range("L11").Value = 10
$ws.Range("M11").Value = 10
$ws.Range("N11").Value = 9
$ws.Range("O11").Value = 6

# Row 12
$ws.Range("B12").Value = 0.04078119335830931
$ws.Range("C12").Value = 0.1801717197739482
$ws.Range("D12").Value = 0.4575852109408234
$ws.Range("E12").Value = 0.06230324688604808
$ws.Range("L12").Value = 10
$ws.Range("M12").Value = 10
$ws.Range("N12").Value = 4
$ws.Range("O12").Value = 3

# Row 13
$ws.Range("B13").Value = 0.1155053710488817
$ws.Range("C13").Value = 0.08243110603500051
$ws.Range("D13").Value = 0.2019409617545528
$ws.Range("E13").Value = 0.1836924325818537
$ws.Range("L13").Value = 10
$ws.Range("M13").Value = 10
$ws.Range("N13").Value = 7
$ws.Range("O13").Value = 5

# Row 14
$ws.Range("B14").Value = 0.0152479675680171
$ws.Range("C14").Value = -0.01308173958600202
$ws.Range("D14").Value = 0.1787191278060382
$ws.Range("E14").Value = 0.1741677061458269
$ws.Range("L14").Value = 10
$ws.Range("M14").Value = 10
$ws.Range("N14").Value = 5
$ws.Range("O14").Value = 4

# Row 15
$ws.Range("B15").Value = -0.04172071846149979
$ws.Range("C15").Value = 0.001014223157493141
$ws.Range("D15").Value = 0.0008212656316868618
$ws.Range("E15").Value = -0.004623178710214178
$ws.Range("L15").Value = 10
$ws.Range("M15").Value = 10
$ws.Range("N15").Value = 6
$ws.Range("O15").Value = 4

# Row 16
$ws.Range("B16").Value = 0.0629556080980055
$ws.Range("C16").Value = -0.001180259424711383
$ws.Range("D16").Value = -0.1067852262153993
$ws.Range("L16").Value = 10
$ws.Range("M16").Value = 10
$ws.Range("N16").Value = 10

# Row 17
$ws.Range("B17").Value = 0.06204703187343687
$ws.Range("C17").Value = 0.04502700709106415
$ws.Range("D17").Value = -0.07255200997378244
$ws.Range("L17").Value = 10
$ws.Range("M17").Value = 10
$ws.Range("N17").Value = 6

